$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Cell Values" sheet: add row 8 with #DIV/0! error values / labels
# ---------------------------------------------------------------------------
$wsValues = $wb.Worksheets.Item("Cell Values")

$wsValues.Range("B8").Value = "#DIV/0!"
$wsValues.Range("C8").Value = "#DIV/0!"
$wsValues.Range("D8").Value = "#DIV/0!"
$wsValues.Range("E8").Value = "#DIV/0!"
$wsValues.Range("F8").Value = "'#DIV/0!"
$wsValues.Range("G8").Value = "'#DIV/0!"

# ---------------------------------------------------------------------------
# 2) Add new "Errors" worksheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsErrors = $wb.Worksheets.Add($null, $lastSheet)
$wsErrors.Name = "Errors"

$wsErrors.Range("B2").Value = "Error value"
$wsErrors.Range("C2").Value = "Formula error"

$wsErrors.Range("B3").Value = "#REF!"
$wsErrors.Range("C3").Formula = "=#REF!+1"

$wsErrors.Range("B4").Value = "#VALUE!"
$wsErrors.Range("C4").Formula = "=""TRUE""*1"

$wsErrors.Range("B5").Value = "#DIV/0!"
$wsErrors.Range("C5").Formula = "=1/0"

$wsErrors.Range("B6").Value = "#NAME?"
$wsErrors.Range("C6").Formula = "=NONEXISTENT.FUNCTION()"

$wsErrors.Range("B7").Value = "#N/A"
$wsErrors.Range("C7").Formula = "=NA()"

$wsErrors.Range("B8").Value = "#NULL!"
$wsErrors.Range("C8").Formula = "=#NULL!+1"

$wsErrors.Range("B9").Value = "#NUM!"
$wsErrors.Range("C9").Formula = "=#NUM!+1"
